$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 15
$ws.Range("B1").Value = 3.250259637832642
$ws.Range("C1").Value = 1.943016171455383
$ws.Range("D1").Value = 1.605194807052612
$ws.Range("E1").Value = 1.516299843788147
